$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Letter date: "September 19, 2025" -> "September 21, 2025"
# ---------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "September 19, 2025") {
        $p.Range.Text = "September 21, 2025"
        break
    }
}

# ---------------------------------------------------------------------
# 2) Split the letterhead mailing-address line
#    "2541 Greenrock Road, Milpitas CA 95035" into two paragraphs:
#    "2541 Greenrock Road" followed by a new "Milpitas, CA 95035"
#    paragraph (same run/paragraph formatting). Only the first
#    occurrence (the letterhead) is touched - the "PROPERTY ADDRESS"
#    table cell further down keeps its original combined text.
# ---------------------------------------------------------------------
$addrIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "2541 Greenrock Road, Milpitas CA 95035") {
        $addrIdx = $i
        break
    }
}
if ($addrIdx -gt 0) {
    $addrPara = $d.Paragraphs.Item($addrIdx)
    $addrPara.Range.Text = "2541 Greenrock Road"
    $addrPara.Range.InsertParagraphAfter()
    $d.Paragraphs.Item($addrIdx + 1).Range.Text = "Milpitas, CA 95035"
}

# ---------------------------------------------------------------------
# 3) Remove the blank "No Spacing" paragraph that immediately follows
#    "Lees Orchard Association Board of Directors".
# ---------------------------------------------------------------------
$boardIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "Board of Directors") {
        $boardIdx = $i
        break
    }
}
if ($boardIdx -gt 0) {
    $blank = $d.Paragraphs.Item($boardIdx + 1)
    if ($blank.Range.Text.Trim() -eq "") {
        $blank.Range.Delete()
    }
}
